$wb = $excel.ActiveWorkbook

# Updated "want-to-go" counts (column F) and "minimum price" values (column G)
# for the exhibition sheet ("展览") and the combined "全部类型" sheet.
# Both sheets contain the same logical rows; the "全部类型" sheet has one
# extra row inserted above row 13, so its row numbers from row 13 onward are
# shifted by +1 relative to the "展览" sheet.

$sheet1Updates = @{
    2  = @{ F = 1899 }
    3  = @{ F = 28 }
    4  = @{ F = 865 }
    6  = @{ F = 50; G = 40 }
    10 = @{ F = 157 }
    13 = @{ F = 4482 }
    15 = @{ F = 37 }
    16 = @{ F = 491 }
    17 = @{ F = 444 }
    18 = @{ F = 14 }
    20 = @{ F = 1118 }
    21 = @{ F = 2243 }
    23 = @{ F = 62 }
    24 = @{ F = 41 }
    25 = @{ G = 50 }
    26 = @{ F = 2191 }
    27 = @{ F = 85 }
    28 = @{ F = 70 }
    29 = @{ F = 20 }
    30 = @{ F = 155 }
    31 = @{ F = 101 }
    32 = @{ F = 39 }
    33 = @{ F = 220 }
}

$sheet4Updates = @{
    2  = @{ F = 1899 }
    3  = @{ F = 28 }
    4  = @{ F = 865 }
    6  = @{ F = 50; G = 40 }
    10 = @{ F = 157 }
    14 = @{ F = 4482 }
    16 = @{ F = 37 }
    17 = @{ F = 491 }
    18 = @{ F = 444 }
    19 = @{ F = 14 }
    21 = @{ F = 1118 }
    22 = @{ F = 2243 }
    24 = @{ F = 62 }
    25 = @{ F = 41 }
    27 = @{ F = 2191 }
    28 = @{ F = 85 }
    29 = @{ F = 70 }
    30 = @{ F = 20 }
    31 = @{ F = 155 }
    32 = @{ F = 101 }
    33 = @{ F = 39 }
    34 = @{ F = 220 }
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $cols = $sheet1Updates[$row]
    foreach ($col in $cols.Keys) {
        $ws1.Range("$col$row").Value = $cols[$col]
    }
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $cols = $sheet4Updates[$row]
    foreach ($col in $cols.Keys) {
        $ws4.Range("$col$row").Value = $cols[$col]
    }
}
